$wb = $excel.ActiveWorkbook

# --- Step 1: rename sheet tabs (two-phase to avoid name collisions) ---
for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Name = "tmp${i}sheetxyz"
}
$newNames = @("summ0", "summ12", "summ26", "summ17", "summ1", "summ4", "summ3", "summ7", "summ9")
for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Name = $newNames[$i-1]
}

# --- Step 2: for each sheet, delete the Single_Female_Parent row, rename
#     Single_Male_Parent -> Single_Parent, and overwrite all coefficient/p values ---

# sheet1.xml (position 1)
$ws = $wb.Worksheets.Item(1)
$ws.Rows.Item(5).Delete()
$ws.Cells.Item(6,1).Value = "HHType_simp[T.Single_Parent]"
$ws.Cells.Item(2,2).Value = 1.732970134219733
$ws.Cells.Item(2,3).Value = 0.5788933368044945
$ws.Cells.Item(3,2).Value = 0.1402070372986627
$ws.Cells.Item(3,3).Value = 0.680347882820815
$ws.Cells.Item(4,2).Value = -1.760880247254324
$ws.Cells.Item(4,3).Value = [double]"4.253347451356256e-10"
$ws.Cells.Item(5,2).Value = -1.084870630545863
$ws.Cells.Item(5,3).Value = 0.0006733268656353179
$ws.Cells.Item(6,2).Value = -1.052150711247818
$ws.Cells.Item(6,3).Value = 0.008769498341023856
$ws.Cells.Item(7,2).Value = -0.1747714873291102
$ws.Cells.Item(7,3).Value = 0.21155821344098
$ws.Cells.Item(8,2).Value = 0.0003151430381457717
$ws.Cells.Item(8,3).Value = 0.0003147787039639631
$ws.Cells.Item(9,2).Value = 0.01692375788382613
$ws.Cells.Item(9,3).Value = 0.0298902246780118
$ws.Cells.Item(10,2).Value = -0.05712179346791828
$ws.Cells.Item(10,3).Value = 0.7757620078218286
$ws.Cells.Item(11,2).Value = 0.8381326990203983
$ws.Cells.Item(11,3).Value = 0.009737538070343419
$ws.Cells.Item(12,2).Value = 0.5847988946481067
$ws.Cells.Item(12,3).Value = 0.167797083548605
$ws.Cells.Item(13,2).Value = [double]"6.381641639757076e-05"
$ws.Cells.Item(13,3).Value = 0.4744277897511214
$ws.Cells.Item(14,2).Value = [double]"-7.218292228718558e-08"
$ws.Cells.Item(14,3).Value = 0.5156888080733488
$ws.Cells.Item(15,2).Value = 0.07183608942017779
$ws.Cells.Item(15,3).Value = 0.574249172883732
$ws.Cells.Item(16,2).Value = 0.05444069907207688
$ws.Cells.Item(16,3).Value = 0.628967460739904
$ws.Cells.Item(17,2).Value = -1.100936880213008
$ws.Cells.Item(17,3).Value = 0.2202659815276907
$ws.Cells.Item(18,2).Value = -0.03539517138408022
$ws.Cells.Item(18,3).Value = 0.3905873323990291
$ws.Cells.Item(19,2).Value = -0.01457359381740956
$ws.Cells.Item(19,3).Value = 0.2457898924542468
$ws.Cells.Item(20,2).Value = 1.300170389712324
$ws.Cells.Item(20,3).Value = 0.7573303756037328
$ws.Cells.Item(21,2).Value = 4.463006843609498
$ws.Cells.Item(21,3).Value = 0.1902963772626065
$ws.Cells.Item(22,2).Value = -0.5649302936575309
$ws.Cells.Item(22,3).Value = 0.777532007569546

# sheet2.xml (position 2)
$ws = $wb.Worksheets.Item(2)
$ws.Rows.Item(5).Delete()
$ws.Cells.Item(6,1).Value = "HHType_simp[T.Single_Parent]"
$ws.Cells.Item(2,2).Value = -0.3961960899081418
$ws.Cells.Item(2,3).Value = 0.8998466713955238
$ws.Cells.Item(3,2).Value = 0.2430411820679466
$ws.Cells.Item(3,3).Value = 0.473810574767594
$ws.Cells.Item(4,2).Value = -1.627523832166488
$ws.Cells.Item(4,3).Value = [double]"6.971966460541229e-09"
$ws.Cells.Item(5,2).Value = -1.037274248873636
$ws.Cells.Item(5,3).Value = 0.0007328466203509502
$ws.Cells.Item(6,2).Value = -0.7342581977727637
$ws.Cells.Item(6,3).Value = 0.06817978232003598
$ws.Cells.Item(7,2).Value = -0.1568264802736347
$ws.Cells.Item(7,3).Value = 0.2632693857831007
$ws.Cells.Item(8,2).Value = 0.0003253318680839094
$ws.Cells.Item(8,3).Value = 0.0001841895442212907
$ws.Cells.Item(9,2).Value = 0.01287595579995072
$ws.Cells.Item(9,3).Value = 0.09352701736229857
$ws.Cells.Item(10,2).Value = 0.1012759286203281
$ws.Cells.Item(10,3).Value = 0.6072135175503812
$ws.Cells.Item(11,2).Value = 0.7047554151469905
$ws.Cells.Item(11,3).Value = 0.03047123682956792
$ws.Cells.Item(12,2).Value = 0.6971090601581782
$ws.Cells.Item(12,3).Value = 0.09558970651208001
$ws.Cells.Item(13,2).Value = [double]"3.587044859919767e-05"
$ws.Cells.Item(13,3).Value = 0.6807685329909774
$ws.Cells.Item(14,2).Value = [double]"2.429753548054712e-08"
$ws.Cells.Item(14,3).Value = 0.8331635416753843
$ws.Cells.Item(15,2).Value = 0.1832556627052955
$ws.Cells.Item(15,3).Value = 0.189942092496569
$ws.Cells.Item(16,2).Value = 0.1426764309384909
$ws.Cells.Item(16,3).Value = 0.2350632530182755
$ws.Cells.Item(17,2).Value = -0.3904795068072336
$ws.Cells.Item(17,3).Value = 0.6654891127590458
$ws.Cells.Item(18,2).Value = -0.03137714679471462
$ws.Cells.Item(18,3).Value = 0.4513563002241504
$ws.Cells.Item(19,2).Value = -0.009271448259748644
$ws.Cells.Item(19,3).Value = 0.4621901993348864
$ws.Cells.Item(20,2).Value = 0.8937795878021645
$ws.Cells.Item(20,3).Value = 0.8345844510266089
$ws.Cells.Item(21,2).Value = 2.776417705180432
$ws.Cells.Item(21,3).Value = 0.4187640026856427
$ws.Cells.Item(22,2).Value = 0.4573927263275616
$ws.Cells.Item(22,3).Value = 0.8211254898672717

# sheet3.xml (position 3)
$ws = $wb.Worksheets.Item(3)
$ws.Rows.Item(5).Delete()
$ws.Cells.Item(6,1).Value = "HHType_simp[T.Single_Parent]"
$ws.Cells.Item(2,2).Value = 2.737927155579248
$ws.Cells.Item(2,3).Value = 0.3912922025438544
$ws.Cells.Item(3,2).Value = 0.4064370737827495
$ws.Cells.Item(3,3).Value = 0.2358612065794189
$ws.Cells.Item(4,2).Value = -1.762411840283054
$ws.Cells.Item(4,3).Value = [double]"5.090374787710027e-10"
$ws.Cells.Item(5,2).Value = -1.092004323674565
$ws.Cells.Item(5,3).Value = 0.0004699581178318277
$ws.Cells.Item(6,2).Value = -0.7085760391487115
$ws.Cells.Item(6,3).Value = 0.08173430763718222
$ws.Cells.Item(7,2).Value = -0.1835197429266717
$ws.Cells.Item(7,3).Value = 0.2070689049431559
$ws.Cells.Item(8,2).Value = 0.0002849451728888131
$ws.Cells.Item(8,3).Value = 0.00110583264850569
$ws.Cells.Item(9,2).Value = 0.01664527944658942
$ws.Cells.Item(9,3).Value = 0.03001891396348412
$ws.Cells.Item(10,2).Value = -0.03695417059227506
$ws.Cells.Item(10,3).Value = 0.8519847711482565
$ws.Cells.Item(11,2).Value = 0.6516346849281778
$ws.Cells.Item(11,3).Value = 0.04091406167241099
$ws.Cells.Item(12,2).Value = 0.5259144941521061
$ws.Cells.Item(12,3).Value = 0.2029179417255137
$ws.Cells.Item(13,2).Value = [double]"-4.47117601395026e-06"
$ws.Cells.Item(13,3).Value = 0.9582031371985662
$ws.Cells.Item(14,2).Value = [double]"-1.094838747537048e-07"
$ws.Cells.Item(14,3).Value = 0.3169628791899134
$ws.Cells.Item(15,2).Value = 0.06858703837627086
$ws.Cells.Item(15,3).Value = 0.5903785865400597
$ws.Cells.Item(16,2).Value = -0.02546582592640796
$ws.Cells.Item(16,3).Value = 0.8232408370602387
$ws.Cells.Item(17,2).Value = 0.4244484634470852
$ws.Cells.Item(17,3).Value = 0.6300695283023481
$ws.Cells.Item(18,2).Value = -0.07229373881370782
$ws.Cells.Item(18,3).Value = 0.08670091323264559
$ws.Cells.Item(19,2).Value = -0.01574650856717114
$ws.Cells.Item(19,3).Value = 0.213624152488729
$ws.Cells.Item(20,2).Value = 7.129627644875446
$ws.Cells.Item(20,3).Value = 0.1014111071412961
$ws.Cells.Item(21,2).Value = 7.219135815346799
$ws.Cells.Item(21,3).Value = 0.03704749002664947
$ws.Cells.Item(22,2).Value = -2.060867177430486
$ws.Cells.Item(22,3).Value = 0.3107081517757572

# sheet4.xml (position 4)
$ws = $wb.Worksheets.Item(4)
$ws.Rows.Item(5).Delete()
$ws.Cells.Item(6,1).Value = "HHType_simp[T.Single_Parent]"
$ws.Cells.Item(2,2).Value = 0.7712558406925705
$ws.Cells.Item(2,3).Value = 0.8047128146928335
$ws.Cells.Item(3,2).Value = 0.2902020922492691
$ws.Cells.Item(3,3).Value = 0.3908227678807443
$ws.Cells.Item(4,2).Value = -1.822333798669052
$ws.Cells.Item(4,3).Value = [double]"9.524083747437283e-11"
$ws.Cells.Item(5,2).Value = -1.257408882326387
$ws.Cells.Item(5,3).Value = [double]"4.31899834371285e-05"
$ws.Cells.Item(6,2).Value = -0.7320268309215129
$ws.Cells.Item(6,3).Value = 0.0806718835888293
$ws.Cells.Item(7,2).Value = -0.2245185137707731
$ws.Cells.Item(7,3).Value = 0.1063367751243373
$ws.Cells.Item(8,2).Value = 0.000300827559913022
$ws.Cells.Item(8,3).Value = 0.0004183335323603265
$ws.Cells.Item(9,2).Value = 0.01349717950884158
$ws.Cells.Item(9,3).Value = 0.07654603127861745
$ws.Cells.Item(10,2).Value = -0.02703894448956958
$ws.Cells.Item(10,3).Value = 0.8914069423722016
$ws.Cells.Item(11,2).Value = 0.7171209759352889
$ws.Cells.Item(11,3).Value = 0.02460125864835196
$ws.Cells.Item(12,2).Value = 0.5644807384429631
$ws.Cells.Item(12,3).Value = 0.1696646890161244
$ws.Cells.Item(13,2).Value = [double]"3.137500732836972e-05"
$ws.Cells.Item(13,3).Value = 0.7175267394180305
$ws.Cells.Item(14,2).Value = [double]"-2.698641685270464e-08"
$ws.Cells.Item(14,3).Value = 0.8164584018388189
$ws.Cells.Item(15,2).Value = 0.1503495011441026
$ws.Cells.Item(15,3).Value = 0.293922982949462
$ws.Cells.Item(16,2).Value = 0.1054858764879442
$ws.Cells.Item(16,3).Value = 0.3893606819872375
$ws.Cells.Item(17,2).Value = -0.7655243895243062
$ws.Cells.Item(17,3).Value = 0.3898785820956613
$ws.Cells.Item(18,2).Value = -0.02175547874184869
$ws.Cells.Item(18,3).Value = 0.596193111434306
$ws.Cells.Item(19,2).Value = -0.009101536911021363
$ws.Cells.Item(19,3).Value = 0.4658664045526016
$ws.Cells.Item(20,2).Value = -0.1848731826067626
$ws.Cells.Item(20,3).Value = 0.9653811074584109
$ws.Cells.Item(21,2).Value = 2.904276307625861
$ws.Cells.Item(21,3).Value = 0.3857115535680703
$ws.Cells.Item(22,2).Value = -0.09299255474077947
$ws.Cells.Item(22,3).Value = 0.962688237199157

# sheet5.xml (position 5)
$ws = $wb.Worksheets.Item(5)
$ws.Rows.Item(5).Delete()
$ws.Cells.Item(6,1).Value = "HHType_simp[T.Single_Parent]"
$ws.Cells.Item(2,2).Value = 2.394734105871072
$ws.Cells.Item(2,3).Value = 0.4491944036095133
$ws.Cells.Item(3,2).Value = 0.2436441284872186
$ws.Cells.Item(3,3).Value = 0.4642984188551967
$ws.Cells.Item(4,2).Value = -1.398370776241264
$ws.Cells.Item(4,3).Value = [double]"5.952077603948829e-07"
$ws.Cells.Item(5,2).Value = -1.019124194538082
$ws.Cells.Item(5,3).Value = 0.001176459034854917
$ws.Cells.Item(6,2).Value = -0.8245068670057483
$ws.Cells.Item(6,3).Value = 0.03595353499726828
$ws.Cells.Item(7,2).Value = -0.1340175983710166
$ws.Cells.Item(7,3).Value = 0.3346298286016055
$ws.Cells.Item(8,2).Value = 0.0003101040288318259
$ws.Cells.Item(8,3).Value = 0.0003178738478414198
$ws.Cells.Item(9,2).Value = 0.01355462091155064
$ws.Cells.Item(9,3).Value = 0.0741633942785749
$ws.Cells.Item(10,2).Value = 0.01195416933023935
$ws.Cells.Item(10,3).Value = 0.9515218016197277
$ws.Cells.Item(11,2).Value = 0.8034573851590593
$ws.Cells.Item(11,3).Value = 0.01298408345374145
$ws.Cells.Item(12,2).Value = 0.67179442531111
$ws.Cells.Item(12,3).Value = 0.1080952834998198
$ws.Cells.Item(13,2).Value = [double]"3.3233263829594e-05"
$ws.Cells.Item(13,3).Value = 0.7009968172175263
$ws.Cells.Item(14,2).Value = [double]"-5.096219138029095e-08"
$ws.Cells.Item(14,3).Value = 0.6478499892786252
$ws.Cells.Item(15,2).Value = 0.09310742283446484
$ws.Cells.Item(15,3).Value = 0.4837034726144441
$ws.Cells.Item(16,2).Value = -0.001923487055754845
$ws.Cells.Item(16,3).Value = 0.9866879500593853
$ws.Cells.Item(17,2).Value = -0.6959610793597378
$ws.Cells.Item(17,3).Value = 0.4281438520467381
$ws.Cells.Item(18,2).Value = -0.06340338861009188
$ws.Cells.Item(18,3).Value = 0.1247398059949281
$ws.Cells.Item(19,2).Value = -0.01374574854703665
$ws.Cells.Item(19,3).Value = 0.2820582873922289
$ws.Cells.Item(20,2).Value = 6.190282064361679
$ws.Cells.Item(20,3).Value = 0.1484003244648255
$ws.Cells.Item(21,2).Value = 7.264296169643756
$ws.Cells.Item(21,3).Value = 0.03417933216405786
$ws.Cells.Item(22,2).Value = -2.845728476894707
$ws.Cells.Item(22,3).Value = 0.1581767004520884

# sheet6.xml (position 6)
$ws = $wb.Worksheets.Item(6)
$ws.Rows.Item(5).Delete()
$ws.Cells.Item(6,1).Value = "HHType_simp[T.Single_Parent]"
$ws.Cells.Item(2,2).Value = 1.796897747652966
$ws.Cells.Item(2,3).Value = 0.5613580097464463
$ws.Cells.Item(3,2).Value = 0.6910942183160805
$ws.Cells.Item(3,3).Value = 0.05343038614788616
$ws.Cells.Item(4,2).Value = -1.763556466156827
$ws.Cells.Item(4,3).Value = [double]"3.260134287542627e-10"
$ws.Cells.Item(5,2).Value = -1.136942425081
$ws.Cells.Item(5,3).Value = 0.0002125068105909029
$ws.Cells.Item(6,2).Value = -0.6588887991465638
$ws.Cells.Item(6,3).Value = 0.105343167895302
$ws.Cells.Item(7,2).Value = -0.3378949322042243
$ws.Cells.Item(7,3).Value = 0.0223353530508004
$ws.Cells.Item(8,2).Value = 0.0003449737849856201
$ws.Cells.Item(8,3).Value = [double]"8.164559315289287e-05"
$ws.Cells.Item(9,2).Value = 0.01725957571818092
$ws.Cells.Item(9,3).Value = 0.02408046954436792
$ws.Cells.Item(10,2).Value = -0.05219585610696161
$ws.Cells.Item(10,3).Value = 0.789425598918057
$ws.Cells.Item(11,2).Value = 0.6427989487262054
$ws.Cells.Item(11,3).Value = 0.04839273346335304
$ws.Cells.Item(12,2).Value = 0.577858437873629
$ws.Cells.Item(12,3).Value = 0.1689907302781245
$ws.Cells.Item(13,2).Value = [double]"7.526324557199427e-06"
$ws.Cells.Item(13,3).Value = 0.9307951482361974
$ws.Cells.Item(14,2).Value = [double]"-8.413019851777506e-08"
$ws.Cells.Item(14,3).Value = 0.441800246374629
$ws.Cells.Item(15,2).Value = 0.08716119643322959
$ws.Cells.Item(15,3).Value = 0.4975780472676816
$ws.Cells.Item(16,2).Value = 0.03376465740054923
$ws.Cells.Item(16,3).Value = 0.7652710138432599
$ws.Cells.Item(17,2).Value = -0.003583840635141982
$ws.Cells.Item(17,3).Value = 0.9967332703955085
$ws.Cells.Item(18,2).Value = -0.04631104088765971
$ws.Cells.Item(18,3).Value = 0.2581323006532701
$ws.Cells.Item(19,2).Value = -0.01213519679694047
$ws.Cells.Item(19,3).Value = 0.3302200385952783
$ws.Cells.Item(20,2).Value = 4.283667472652061
$ws.Cells.Item(20,3).Value = 0.3031463851406742
$ws.Cells.Item(21,2).Value = 5.814519151576604
$ws.Cells.Item(21,3).Value = 0.08140375282079484
$ws.Cells.Item(22,2).Value = -1.623482284984084
$ws.Cells.Item(22,3).Value = 0.4035992030458452

# sheet7.xml (position 7)
$ws = $wb.Worksheets.Item(7)
$ws.Rows.Item(5).Delete()
$ws.Cells.Item(6,1).Value = "HHType_simp[T.Single_Parent]"
$ws.Cells.Item(2,2).Value = 2.700809469212552
$ws.Cells.Item(2,3).Value = 0.3792386669367031
$ws.Cells.Item(3,2).Value = 0.2334202282954864
$ws.Cells.Item(3,3).Value = 0.4813907732449674
$ws.Cells.Item(4,2).Value = -1.538315289985376
$ws.Cells.Item(4,3).Value = [double]"3.582559869724834e-08"
$ws.Cells.Item(5,2).Value = -0.9659674870167219
$ws.Cells.Item(5,3).Value = 0.001627827548398885
$ws.Cells.Item(6,2).Value = -0.7163771066396583
$ws.Cells.Item(6,3).Value = 0.06776077617013587
$ws.Cells.Item(7,2).Value = -0.1192715068815926
$ws.Cells.Item(7,3).Value = 0.387656860589067
$ws.Cells.Item(8,2).Value = 0.0003376199115320621
$ws.Cells.Item(8,3).Value = [double]"7.617628087165788e-05"
$ws.Cells.Item(9,2).Value = 0.01219878116216787
$ws.Cells.Item(9,3).Value = 0.1044518479806271
$ws.Cells.Item(10,2).Value = -0.1501805045730079
$ws.Cells.Item(10,3).Value = 0.4409304276669024
$ws.Cells.Item(11,2).Value = 0.4911012580170867
$ws.Cells.Item(11,3).Value = 0.1211472973633454
$ws.Cells.Item(12,2).Value = 0.5731365879516369
$ws.Cells.Item(12,3).Value = 0.1612802772758984
$ws.Cells.Item(13,2).Value = [double]"3.361049528881585e-05"
$ws.Cells.Item(13,3).Value = 0.6960362052254101
$ws.Cells.Item(14,2).Value = [double]"-1.056556842092603e-07"
$ws.Cells.Item(14,3).Value = 0.3413621485071723
$ws.Cells.Item(15,2).Value = 0.0531293481969306
$ws.Cells.Item(15,3).Value = 0.6810002603326446
$ws.Cells.Item(16,2).Value = 0.01004058836155924
$ws.Cells.Item(16,3).Value = 0.9288566818143722
$ws.Cells.Item(17,2).Value = -0.6488118569682828
$ws.Cells.Item(17,3).Value = 0.4486506420033632
$ws.Cells.Item(18,2).Value = -0.04700597779239762
$ws.Cells.Item(18,3).Value = 0.2437868851383325
$ws.Cells.Item(19,2).Value = -0.0150682555996215
$ws.Cells.Item(19,3).Value = 0.223420523819406
$ws.Cells.Item(20,2).Value = 3.703655746381877
$ws.Cells.Item(20,3).Value = 0.3735790509902489
$ws.Cells.Item(21,2).Value = 5.928323918181067
$ws.Cells.Item(21,3).Value = 0.07660641023453144
$ws.Cells.Item(22,2).Value = -1.783580042546859
$ws.Cells.Item(22,3).Value = 0.3646161455179876

# sheet8.xml (position 8)
$ws = $wb.Worksheets.Item(8)
$ws.Rows.Item(5).Delete()
$ws.Cells.Item(6,1).Value = "HHType_simp[T.Single_Parent]"
$ws.Cells.Item(2,2).Value = 2.386227134354089
$ws.Cells.Item(2,3).Value = 0.4414850706890165
$ws.Cells.Item(3,2).Value = 0.2998033395361526
$ws.Cells.Item(3,3).Value = 0.3767419780582052
$ws.Cells.Item(4,2).Value = -1.767063615186987
$ws.Cells.Item(4,3).Value = [double]"4.43989856001058e-10"
$ws.Cells.Item(5,2).Value = -1.117137083825381
$ws.Cells.Item(5,3).Value = 0.0003193410968254279
$ws.Cells.Item(6,2).Value = -0.7275600429451156
$ws.Cells.Item(6,3).Value = 0.07045744137703795
$ws.Cells.Item(7,2).Value = -0.2269271953349747
$ws.Cells.Item(7,3).Value = 0.1102300102203191
$ws.Cells.Item(8,2).Value = 0.000287057388653829
$ws.Cells.Item(8,3).Value = 0.0007250506970414987
$ws.Cells.Item(9,2).Value = 0.01370085244762884
$ws.Cells.Item(9,3).Value = 0.07051372207725605
$ws.Cells.Item(10,2).Value = -0.1317574219416145
$ws.Cells.Item(10,3).Value = 0.5039222015281231
$ws.Cells.Item(11,2).Value = 0.6941291542466532
$ws.Cells.Item(11,3).Value = 0.03323369767344516
$ws.Cells.Item(12,2).Value = 0.6049175199331389
$ws.Cells.Item(12,3).Value = 0.1478516355192141
$ws.Cells.Item(13,2).Value = [double]"2.587198772067267e-05"
$ws.Cells.Item(13,3).Value = 0.7644252046002229
$ws.Cells.Item(14,2).Value = [double]"-1.359031837116362e-07"
$ws.Cells.Item(14,3).Value = 0.2120813720014804
$ws.Cells.Item(15,2).Value = 0.0617216224870991
$ws.Cells.Item(15,3).Value = 0.6272721880968373
$ws.Cells.Item(16,2).Value = 0.02599967551708489
$ws.Cells.Item(16,3).Value = 0.8157734030698627
$ws.Cells.Item(17,2).Value = -0.9623712507156256
$ws.Cells.Item(17,3).Value = 0.2633880333852683
$ws.Cells.Item(18,2).Value = -0.0382483258912747
$ws.Cells.Item(18,3).Value = 0.3473874394737784
$ws.Cells.Item(19,2).Value = -0.01255807817947547
$ws.Cells.Item(19,3).Value = 0.3136734410721118
$ws.Cells.Item(20,2).Value = 2.734076062030778
$ws.Cells.Item(20,3).Value = 0.5138896951289413
$ws.Cells.Item(21,2).Value = 5.093697591283347
$ws.Cells.Item(21,3).Value = 0.1289511753880163
$ws.Cells.Item(22,2).Value = -1.218656063961704
$ws.Cells.Item(22,3).Value = 0.538479948068596

# sheet9.xml (position 9)
$ws = $wb.Worksheets.Item(9)
$ws.Rows.Item(5).Delete()
$ws.Cells.Item(6,1).Value = "HHType_simp[T.Single_Parent]"
$ws.Cells.Item(2,2).Value = 1.26914307491598
$ws.Cells.Item(2,3).Value = 0.6821519595807368
$ws.Cells.Item(3,2).Value = 0.01681133396938332
$ws.Cells.Item(3,3).Value = 0.960093721146796
$ws.Cells.Item(4,2).Value = -1.626513472627707
$ws.Cells.Item(4,3).Value = [double]"1.217615661987365e-08"
$ws.Cells.Item(5,2).Value = -0.8327592292171645
$ws.Cells.Item(5,3).Value = 0.007864065757980819
$ws.Cells.Item(6,2).Value = -0.7169046238032213
$ws.Cells.Item(6,3).Value = 0.07399157808507983
$ws.Cells.Item(7,2).Value = -0.02038795964182197
$ws.Cells.Item(7,3).Value = 0.8883056469267318
$ws.Cells.Item(8,2).Value = 0.0002815807125481109
$ws.Cells.Item(8,3).Value = 0.001502870894809775
$ws.Cells.Item(9,2).Value = 0.0165945000024058
$ws.Cells.Item(9,3).Value = 0.03359212234493894
$ws.Cells.Item(10,2).Value = -0.08939833500704152
$ws.Cells.Item(10,3).Value = 0.6493659123560596
$ws.Cells.Item(11,2).Value = 0.6958198561356732
$ws.Cells.Item(11,3).Value = 0.03230384163859094
$ws.Cells.Item(12,2).Value = 0.3968436294878626
$ws.Cells.Item(12,3).Value = 0.3452649407646506
$ws.Cells.Item(13,2).Value = [double]"6.121152333256449e-05"
$ws.Cells.Item(13,3).Value = 0.4824254387305654
$ws.Cells.Item(14,2).Value = [double]"-6.777484229395552e-08"
$ws.Cells.Item(14,3).Value = 0.5368560375722731
$ws.Cells.Item(15,2).Value = 0.06954535249614767
$ws.Cells.Item(15,3).Value = 0.5862235277596806
$ws.Cells.Item(16,2).Value = 0.04472930059457696
$ws.Cells.Item(16,3).Value = 0.689579653151404
$ws.Cells.Item(17,2).Value = -0.7044125204408943
$ws.Cells.Item(17,3).Value = 0.4246849333618139
$ws.Cells.Item(18,2).Value = -0.04672968643561151
$ws.Cells.Item(18,3).Value = 0.2492285016830645
$ws.Cells.Item(19,2).Value = -0.01171133437485694
$ws.Cells.Item(19,3).Value = 0.3451439612436517
$ws.Cells.Item(20,2).Value = 3.016535889409668
$ws.Cells.Item(20,3).Value = 0.4662202478522907
$ws.Cells.Item(21,2).Value = 5.332276011594635
$ws.Cells.Item(21,3).Value = 0.112585519984966
$ws.Cells.Item(22,2).Value = -1.028132853197482
$ws.Cells.Item(22,3).Value = 0.6029646484566316
